$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.6684293746948242
$ws.Range("D3").Value = 0.2649017572402954
$ws.Range("D4").Value = 0.845078706741333
$ws.Range("D5").Value = 0.6489423513412476
$ws.Range("D6").Value = 0.06797409057617188
$ws.Range("D7").Value = 0.9090899229049683
$ws.Range("D8").Value = 0.7793173789978027
$ws.Range("D9").Value = 0.4656662344932556
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 0.3958871364593506
$ws.Range("D11").Value = 0.2400163263082504
$ws.Range("D12").Value = 0.1454481184482574
$ws.Range("D13").Value = 0.8413236141204834
$ws.Range("D14").Value = 0.9368212223052979
$ws.Range("C15").Value = 0
$ws.Range("D15").Value = 0.4560641944408417
$ws.Range("D16").Value = 0.2331621944904327
$ws.Range("D17").Value = 0.805791974067688
$ws.Range("C18").Value = 0
$ws.Range("D18").Value = 0.3247884809970856
$ws.Range("C19").Value = 1
$ws.Range("D19").Value = 0.706403374671936
$ws.Range("D20").Value = 0.4411842823028564
$ws.Range("D21").Value = 0.1344810873270035
$ws.Range("D22").Value = 0.5167005062103271
$ws.Range("D23").Value = 0.8343610763549805
$ws.Range("D24").Value = 0.4344661235809326
$ws.Range("D25").Value = 0.830958366394043
$ws.Range("D26").Value = 0.8979665637016296
$ws.Range("D27").Value = 0.2046066969633102
$ws.Range("D28").Value = 0.7010354399681091
$ws.Range("D29").Value = 0.1762720197439194
$ws.Range("D30").Value = 0.8072061538696289
$ws.Range("D31").Value = 0.1294015794992447
$ws.Range("D32").Value = 0.6676579713821411
$ws.Range("D33").Value = 0.6594008207321167
$ws.Range("D34").Value = 0.9113715291023254
$ws.Range("D35").Value = 0.8335343599319458
$ws.Range("D36").Value = 0.2622500360012054
$ws.Range("D37").Value = 0.6606177091598511
$ws.Range("D38").Value = 0.9014912843704224
$ws.Range("D39").Value = 0.9648342728614807
$ws.Range("D40").Value = 0.8533512353897095
$ws.Range("D41").Value = 0.9497278332710266
$ws.Range("C42").Value = 0
$ws.Range("D42").Value = 0.4550877809524536
$ws.Range("D43").Value = 0.7018100619316101
$ws.Range("D44").Value = 0.7118388414382935
$ws.Range("D45").Value = 0.146382674574852
$ws.Range("D46").Value = 0.8913187980651855
$ws.Range("D47").Value = 0.8311642408370972
$ws.Range("D48").Value = 0.622575044631958
$ws.Range("D49").Value = 0.3789398372173309
$ws.Range("D50").Value = 0.2491694092750549
$ws.Range("D51").Value = 0.5190890431404114
$ws.Range("D52").Value = 0.1144847124814987
$ws.Range("D53").Value = 0.599481999874115
$ws.Range("D54").Value = 0.9440380334854126
$ws.Range("D55").Value = 0.2743572890758514
$ws.Range("C56").Value = 0
$ws.Range("D56").Value = 0.4821081757545471
$ws.Range("D57").Value = 0.2576489150524139
$ws.Range("D58").Value = 0.6313918232917786
$ws.Range("D59").Value = 0.5965096354484558
$ws.Range("D60").Value = 0.9467431306838989
$ws.Range("D61").Value = 0.6300801038742065
$ws.Range("D62").Value = 0.6403828859329224
$ws.Range("D63").Value = 0.688056468963623
$ws.Range("D64").Value = 0.1478163003921509
$ws.Range("D65").Value = 0.1426456719636917
$ws.Range("D66").Value = 0.7692208290100098
$ws.Range("D67").Value = 0.4048758149147034
$ws.Range("D68").Value = 0.4813310205936432
$ws.Range("D69").Value = 0.7104834318161011
$ws.Range("D70").Value = 0.2862207293510437
$ws.Range("D71").Value = 0.3769761025905609
$ws.Range("D72").Value = 0.4086599946022034
$ws.Range("D73").Value = 0.8203609585762024
$ws.Range("D74").Value = 0.4387290179729462
$ws.Range("D75").Value = 0.1037182807922363
$ws.Range("D76").Value = 0.3555331528186798
$ws.Range("C77").Value = 0
$ws.Range("D77").Value = 0.3520670831203461
$ws.Range("D78").Value = 0.7589368224143982
$ws.Range("D79").Value = 0.1492019593715668
